# Chokhatauri.xlsx — "upgrade left table until javakheti"
#
# 1. Rename the only worksheet from the placeholder "1" to "Chokhatauri".
# 2. The "Urban" row (row 6) becomes entirely confidential/unavailable data,
#    so every year column (D:O) is set to the "…" marker.
# 3. In the "Rural" row (row 7) the last two observed years (J7:K7 -> 2018,
#    2019) are also masked to "…", while the remaining numbers are left as-is.
# 4. The blank spacer row (row 8) is removed so the "Note:" row moves up
#    from row 9 to row 8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Chokhatauri"

$ws.Range("D6:O6").Value = "…"
$ws.Range("J7:K7").Value = "…"

$ws.Rows("8").Delete()
